$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the s_vals sheet, mirroring the existing
# "sum" column (G): same header styling, and a numeric 0 for each data row.

# Copy G1's formatting (bold font, border, centered/top alignment) onto H1
# so the new header cell reuses the existing header style.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
